$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$cs = $s.ThemeColorScheme
$cs.Item(3).RGB = 6908522
Write-Output ("new color3=" + $cs.Item(3).RGB)
